$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45, shifting existing rows 45-116 down to 46-117.
$ws.Cells.Item(45, 4).EntireRow.Insert()

# Populate the newly inserted row 45 with its data.
$ws.Cells.Item(45, 1).Value = 11
$ws.Cells.Item(45, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(45, 3).Value = "Bíobío"
$ws.Cells.Item(45, 4).Value = 44540
$ws.Cells.Item(45, 5).Value = 8
$ws.Cells.Item(45, 6).Value = 100112003
$ws.Cells.Item(45, 7).Value = "Ajo"
$ws.Cells.Item(45, 8).Value = "Chino"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 180
$ws.Cells.Item(45, 11).Value = 16000
$ws.Cells.Item(45, 12).Value = 17000
$ws.Cells.Item(45, 13).Value = 16556
$ws.Cells.Item(45, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(45, 15).Value = "China"
$ws.Cells.Item(45, 16).Value = 1656
$ws.Cells.Item(45, 17).Value = 10
$ws.Cells.Item(45, 18).Value = "Hortaliza"
